# AIE/F/14-IIR: Updated for 2023.1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicated closing parenthesis in the header text.
$ws.Range("A9").Value = "API Throughput (Msa/sec)"

# Update measured cycle counts (row 7) for the new 2023.1 results.
$ws.Range("B7").Value = 187
$ws.Range("C7").Value = 492
$ws.Range("D7").Value = 940
$ws.Range("E7").Value = 1836
$ws.Range("F7").Value = 3628
$ws.Range("G7").Value = 7212
$ws.Range("H7").Value = 14379

# Move the active selection to A14, matching the updated view state.
$ws.Range("A14").Select()
